# Adding data for 2021 year + modifications
#
# 1) Set the selection on the existing "data2020" sheet to C6 (it loses the
#    "active tab" status once the new sheet becomes active).
# 2) Insert a brand-new "data2021" worksheet as the first tab and populate
#    it with the 2021 pollution data (mirrors the layout of the other
#    yearly sheets, but the header of column A reads "niewlasciwe pole").

$wb = $excel.ActiveWorkbook

# --- Step 1: fix up the selection on the data2020 sheet -------------------
$ws2020 = $wb.Worksheets.Item("data2020")
$ws2020.Activate()
$ws2020.Range("C6").Select()

# --- Step 2: create the new data2021 sheet as the first tab ---------------
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "data2021"
$newSheet.Move($wb.Worksheets.Item(1))

$data = @(
    @('niewlasciwe pole', 'area', 'co2', 'metan', 'n2o', 'so2', 'no', 'co'),
    @('dolnoslaskie', 19947, 24306.62, 42.2, 3.18, 26.56, 37.65, 158.11000000000001),
    @('kujawsko-pomorskie', 17972, 15750.88, 62.3, 6.2, 27.18, 44.77, 175.51),
    @('lubelskie', 25122, 12270.73, 130.4, 6.11, 19.91, 31.41, 145.78),
    @('lubuskie', 13988, 4986.79, 40.200000000000003, 1.68, 3.98, 12.54, 61.75),
    @('lodzkie', 18219, 48236.01, 93.93, 5.26, 67.25, 66.62, 229.63),
    @('malopolskie', 15183, 17924.34, 45.04, 2.35, 22.78, 32.44, 156.94999999999999),
    @('mazowieckie', 35558, 58045.52, 148.77000000000001, 9.9, 49.15, 87.88, 229.28),
    @('opolskie', 9412, 20797.11, 23.12, 2.76, 13.43, 25.71, 68.319999999999993),
    @('podkarpackie', 17846, 8028.31, 31.42, 1.54, 16.39, 19.260000000000002, 242.66),
    @('podlaskie', 20187, 4119.51, 104.07, 7.32, 7.16, 17.600000000000001, 68.64),
    @('pomorskie', 18310, 13564.65, 53.3, 3.92, 23.21, 35, 161.24),
    @('slaskie', 12333, 50417.5, 488.88, 2.61, 44.31, 55.84, 219.66),
    @('swietokrzyskie', 11711, 15333.41, 22.97, 1.87, 16.02, 22.2, 68.98),
    @('warminsko-mazurskie', 24173, 5524.27, 58.35, 5.67, 12.09, 21.25, 124.8),
    @('wielkopolskie', 29826, 21070.68, 141.33000000000001, 11.73, 23.5, 50.53, 213.79),
    @('zachodniopomorskie', 22892, 11199.94, 38.89, 5.69, 19.47, 30.68, 195.65)
)

for ($r = 0; $r -lt $data.Count; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Count; $c++) {
        $newSheet.Cells.Item($r + 1, $c + 1).Value = $row[$c]
    }
}
